$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Misc. Data" column header in J1
$ws.Range("J1").Value = 'Misc. Data'

# Row 2 (Anthony V Das et al.) - Authors field (E2) now uses the final,
# correctly-padded author string; new Misc. Data cell (J2) is empty
$ws.Range("E2").Value = '[Anthony V%Das%NULL%1,                         Padmaja K%Rani%NULL%1,                         Pravin K%Vaddavalli%NULL%1]'
$ws.Range("J2").Value = ''

# Row 3 (Gagan Kalra et al.) - Authors field (E3) now uses the final,
# correctly-padded author string; new Misc. Data cell (J3) is empty
$ws.Range("E3").Value = '[Gagan%Kalra%NULL%1,                         Andrew M.%Williams%NULL%1,                         Patrick W.%Commiskey%NULL%1,                         Eve M. R.%Bowers%NULL%1,                         Tadhg%Schempf%NULL%1,                         José-Alain%Sahel%NULL%1,                         Evan L.%Waxman%waxmane@upmc.edu%1,                         Roxana%Fu%fur3@upmc.edu%1]'
$ws.Range("J3").Value = ''
